$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-0.244, 0.053, 0.554, -170, -10, 4),
    @(-0.206, -0.368, -0.067, 97, 13, 72),
    @(0.594, 0.03, -0.209, 105, -25, 36),
    @(-0.177, 0.145, 0.001, 155, -14, -176),
    @(0.476, 0.203, -0.188, -179, 2, -150),
    @(-0.037, 0.626, 0.24, -23, 87, 109),
    @(-0.081, 0.199, -0.201, 95, 29, 101),
    @(0.042, -0.083, 0.06, 142, 6, -164),
    @(0.054, -0.15, -0.095, 144, -10, 103),
    @(0.299, 0.5659999999999999, 0.355, 54, 1, 89)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}
